$wb = $excel.ActiveWorkbook

# ---- Sheet: Overview ----
$ws = $wb.Worksheets.Item("Overview")

# Update cell values for rows 4-7
$ws.Range("A4").Value = '8a05a72f-0702-477e-92c2-46140fea9d0a.md'
$ws.Range("B4").Value = 'Ready for handoff'
$ws.Range("C4").Value = 'Ready for handoff'
$ws.Range("D4").Value = '2016-17-17 22:17:54'

$ws.Range("A5").Value = 'b8cc38f3-8078-488f-a3d7-a5ee88996c60.md'
$ws.Range("B5").Value = 'Handed back: in sync with en-US'
$ws.Range("C5").Value = 'Handed back: in sync with en-US'
$ws.Range("D5").Value = '2016-16-17 22:16:19'

$ws.Range("A6").Value = '20f88ce9-2e40-40ac-af6e-41c0a53aab0f.md'
$ws.Range("B6").Value = 'Ready for handoff'
$ws.Range("C6").Value = 'Ready for handoff'
$ws.Range("D6").Value = '2016-17-17 22:17:54'

$ws.Range("A7").Value = '5704218f-4b11-4c3e-ae00-dedbfae3387a.md'
$ws.Range("B7").Value = 'Ready for handoff'
$ws.Range("C7").Value = 'Ready for handoff'
$ws.Range("D7").Value = '2016-17-17 22:17:54'

# Update hyperlink display text for rows 4-7
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    switch ($addr) {
        '$A$5' { $hl.TextToDisplay = 'b8cc38f3-8078-488f-a3d7-a5ee88996c60.md' }
        '$A$6' { $hl.TextToDisplay = '20f88ce9-2e40-40ac-af6e-41c0a53aab0f.md' }
        '$A$7' { $hl.TextToDisplay = '5704218f-4b11-4c3e-ae00-dedbfae3387a.md' }
        '$A$4' { $hl.TextToDisplay = '8a05a72f-0702-477e-92c2-46140fea9d0a.md' }
    }
}

# ---- Sheet: zh-cn ----
$ws = $wb.Worksheets.Item("zh-cn")

# Update cell values for rows 4-7
$ws.Range("A4").Value = '8a05a72f-0702-477e-92c2-46140fea9d0a.md'
$ws.Range("B4").Value = '.md'
$ws.Range("C4").Value = 'Ready for handoff'
$ws.Range("D4").Value = '8a05a72f-0702-477e-92c2-46140fea9d0a.f9a574f5491ae4f28d789074703643d8efad514d.zh-cn.xlf'
$ws.Range("E4").Value = '2016-03-17 22:17:50'
$ws.Range("F4").Value = '8a05a72f-0702-477e-92c2-46140fea9d0a.md'
$ws.Range("G4").Value = '8a05a72f-0702-477e-92c2-46140fea9d0a.f9a574f5491ae4f28d789074703643d8efad514d.zh-cn.xlf'
$ws.Range("H4").Value = '2016-03-17 22:18:10'
$ws.Range("I4").Value = 'Include'

$ws.Range("A5").Value = 'b8cc38f3-8078-488f-a3d7-a5ee88996c60.md'
$ws.Range("B5").Value = '.md'
$ws.Range("C5").Value = 'Handed back: in sync with en-US'
$ws.Range("D5").Value = 'b8cc38f3-8078-488f-a3d7-a5ee88996c60.d3541b304718d0d615f34f432413656b28b6736b.zh-cn.xlf'
$ws.Range("E5").Value = '2016-03-17 22:16:16'
$ws.Range("F5").Value = 'b8cc38f3-8078-488f-a3d7-a5ee88996c60.md'
$ws.Range("G5").Value = 'b8cc38f3-8078-488f-a3d7-a5ee88996c60.d3541b304718d0d615f34f432413656b28b6736b.zh-cn.xlf'
$ws.Range("H5").Value = '2016-03-17 22:16:34'
$ws.Range("I5").Value = 'Include'

$ws.Range("A6").Value = '20f88ce9-2e40-40ac-af6e-41c0a53aab0f.md'
$ws.Range("B6").Value = '.md'
$ws.Range("C6").Value = 'Ready for handoff'
$ws.Range("D6").Value = '20f88ce9-2e40-40ac-af6e-41c0a53aab0f.8a9b7658475068434230ce69758c9384275db8d8.zh-cn.xlf'
$ws.Range("E6").Value = '2016-03-17 22:17:50'
$ws.Range("H6").Value = '0001-01-01 00:00:00'
$ws.Range("I6").Value = 'Include'

$ws.Range("A7").Value = '5704218f-4b11-4c3e-ae00-dedbfae3387a.md'
$ws.Range("B7").Value = '.md'
$ws.Range("C7").Value = 'Ready for handoff'
$ws.Range("D7").Value = '5704218f-4b11-4c3e-ae00-dedbfae3387a.99316c6d6d154cc6858782d5d54c030655d47352.zh-cn.xlf'
$ws.Range("E7").Value = '2016-03-17 22:17:50'
$ws.Range("H7").Value = '0001-01-01 00:00:00'
$ws.Range("I7").Value = 'Include'

# Update hyperlink display text for rows 4-7
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    switch ($addr) {
        '$A$5' { $hl.TextToDisplay = 'b8cc38f3-8078-488f-a3d7-a5ee88996c60.md' }
        '$B$5' { $hl.TextToDisplay = '.md' }
        '$D$5' { $hl.TextToDisplay = 'b8cc38f3-8078-488f-a3d7-a5ee88996c60.d3541b304718d0d615f34f432413656b28b6736b.zh-cn.xlf' }
        '$F$5' { $hl.TextToDisplay = 'b8cc38f3-8078-488f-a3d7-a5ee88996c60.md' }
        '$G$5' { $hl.TextToDisplay = 'b8cc38f3-8078-488f-a3d7-a5ee88996c60.d3541b304718d0d615f34f432413656b28b6736b.zh-cn.xlf' }
        '$A$6' { $hl.TextToDisplay = '20f88ce9-2e40-40ac-af6e-41c0a53aab0f.md' }
        '$B$6' { $hl.TextToDisplay = '.md' }
        '$D$6' { $hl.TextToDisplay = '20f88ce9-2e40-40ac-af6e-41c0a53aab0f.8a9b7658475068434230ce69758c9384275db8d8.zh-cn.xlf' }
        '$A$7' { $hl.TextToDisplay = '5704218f-4b11-4c3e-ae00-dedbfae3387a.md' }
        '$B$7' { $hl.TextToDisplay = '.md' }
        '$D$7' { $hl.TextToDisplay = '5704218f-4b11-4c3e-ae00-dedbfae3387a.99316c6d6d154cc6858782d5d54c030655d47352.zh-cn.xlf' }
        '$A$4' { $hl.TextToDisplay = '8a05a72f-0702-477e-92c2-46140fea9d0a.md' }
        '$B$4' { $hl.TextToDisplay = '.md' }
        '$D$4' { $hl.TextToDisplay = '8a05a72f-0702-477e-92c2-46140fea9d0a.f9a574f5491ae4f28d789074703643d8efad514d.zh-cn.xlf' }
        '$F$4' { $hl.TextToDisplay = '8a05a72f-0702-477e-92c2-46140fea9d0a.md' }
        '$G$4' { $hl.TextToDisplay = '8a05a72f-0702-477e-92c2-46140fea9d0a.f9a574f5491ae4f28d789074703643d8efad514d.zh-cn.xlf' }
    }
}

# ---- Sheet: de-de ----
$ws = $wb.Worksheets.Item("de-de")

# Update cell values for rows 4-7
$ws.Range("A4").Value = '8a05a72f-0702-477e-92c2-46140fea9d0a.md'
$ws.Range("B4").Value = '.md'
$ws.Range("C4").Value = 'Ready for handoff'
$ws.Range("D4").Value = '8a05a72f-0702-477e-92c2-46140fea9d0a.f9a574f5491ae4f28d789074703643d8efad514d.de-de.xlf'
$ws.Range("E4").Value = '2016-03-17 22:17:54'
$ws.Range("F4").Value = '8a05a72f-0702-477e-92c2-46140fea9d0a.md'
$ws.Range("G4").Value = '8a05a72f-0702-477e-92c2-46140fea9d0a.f9a574f5491ae4f28d789074703643d8efad514d.de-de.xlf'
$ws.Range("H4").Value = '2016-03-17 22:18:16'
$ws.Range("I4").Value = 'Include'

$ws.Range("A5").Value = 'b8cc38f3-8078-488f-a3d7-a5ee88996c60.md'
$ws.Range("B5").Value = '.md'
$ws.Range("C5").Value = 'Handed back: in sync with en-US'
$ws.Range("D5").Value = 'b8cc38f3-8078-488f-a3d7-a5ee88996c60.d3541b304718d0d615f34f432413656b28b6736b.de-de.xlf'
$ws.Range("E5").Value = '2016-03-17 22:16:19'
$ws.Range("F5").Value = 'b8cc38f3-8078-488f-a3d7-a5ee88996c60.md'
$ws.Range("G5").Value = 'b8cc38f3-8078-488f-a3d7-a5ee88996c60.d3541b304718d0d615f34f432413656b28b6736b.de-de.xlf'
$ws.Range("H5").Value = '2016-03-17 22:16:40'
$ws.Range("I5").Value = 'Include'

$ws.Range("A6").Value = '20f88ce9-2e40-40ac-af6e-41c0a53aab0f.md'
$ws.Range("B6").Value = '.md'
$ws.Range("C6").Value = 'Ready for handoff'
$ws.Range("D6").Value = '20f88ce9-2e40-40ac-af6e-41c0a53aab0f.8a9b7658475068434230ce69758c9384275db8d8.de-de.xlf'
$ws.Range("E6").Value = '2016-03-17 22:17:54'
$ws.Range("H6").Value = '0001-01-01 00:00:00'
$ws.Range("I6").Value = 'Include'

$ws.Range("A7").Value = '5704218f-4b11-4c3e-ae00-dedbfae3387a.md'
$ws.Range("B7").Value = '.md'
$ws.Range("C7").Value = 'Ready for handoff'
$ws.Range("D7").Value = '5704218f-4b11-4c3e-ae00-dedbfae3387a.99316c6d6d154cc6858782d5d54c030655d47352.de-de.xlf'
$ws.Range("E7").Value = '2016-03-17 22:17:54'
$ws.Range("H7").Value = '0001-01-01 00:00:00'
$ws.Range("I7").Value = 'Include'

# Update hyperlink display text for rows 4-7
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    switch ($addr) {
        '$A$5' { $hl.TextToDisplay = 'b8cc38f3-8078-488f-a3d7-a5ee88996c60.md' }
        '$B$5' { $hl.TextToDisplay = '.md' }
        '$D$5' { $hl.TextToDisplay = 'b8cc38f3-8078-488f-a3d7-a5ee88996c60.d3541b304718d0d615f34f432413656b28b6736b.de-de.xlf' }
        '$F$5' { $hl.TextToDisplay = 'b8cc38f3-8078-488f-a3d7-a5ee88996c60.md' }
        '$G$5' { $hl.TextToDisplay = 'b8cc38f3-8078-488f-a3d7-a5ee88996c60.d3541b304718d0d615f34f432413656b28b6736b.de-de.xlf' }
        '$A$6' { $hl.TextToDisplay = '20f88ce9-2e40-40ac-af6e-41c0a53aab0f.md' }
        '$B$6' { $hl.TextToDisplay = '.md' }
        '$D$6' { $hl.TextToDisplay = '20f88ce9-2e40-40ac-af6e-41c0a53aab0f.8a9b7658475068434230ce69758c9384275db8d8.de-de.xlf' }
        '$A$7' { $hl.TextToDisplay = '5704218f-4b11-4c3e-ae00-dedbfae3387a.md' }
        '$B$7' { $hl.TextToDisplay = '.md' }
        '$D$7' { $hl.TextToDisplay = '5704218f-4b11-4c3e-ae00-dedbfae3387a.99316c6d6d154cc6858782d5d54c030655d47352.de-de.xlf' }
        '$A$4' { $hl.TextToDisplay = '8a05a72f-0702-477e-92c2-46140fea9d0a.md' }
        '$B$4' { $hl.TextToDisplay = '.md' }
        '$D$4' { $hl.TextToDisplay = '8a05a72f-0702-477e-92c2-46140fea9d0a.f9a574f5491ae4f28d789074703643d8efad514d.de-de.xlf' }
        '$F$4' { $hl.TextToDisplay = '8a05a72f-0702-477e-92c2-46140fea9d0a.md' }
        '$G$4' { $hl.TextToDisplay = '8a05a72f-0702-477e-92c2-46140fea9d0a.f9a574f5491ae4f28d789074703643d8efad514d.de-de.xlf' }
    }
}
